# Update the "Average SAT Score" axis label -> "Admission Rate", and the
# four axis tick labels (700/1700/1200/1000) -> percentages (0%/100%/75%/25%),
# repositioning/resizing the textboxes to match the new, shorter labels.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Conversion helper: EMU -> points (PowerPoint COM works in points, stored
# as single-precision floats). Add a tiny sub-EMU bias so that the
# single-precision round-trip back to EMU lands exactly on the target
# value instead of truncating one EMU short.
function EmuToPt([double]$emu) {
    $bias = 0.5 / 914400 * 72
    return ($emu / 914400 * 72) + $bias
}

# --- "Average SAT Score" -> "Admission Rate" (Rectangle 2, id=3) ---
$title = $s.Shapes.Item(41)
$title.Left = EmuToPt 767332
$title.Top = EmuToPt 3243289
$title.Width = EmuToPt 1800493
$title.Height = EmuToPt 369332
$title.TextFrame.TextRange.Text = "Admission Rate"

# --- "700" -> "0%" (TextBox 96, id=97) ---
$lbl700 = $s.Shapes.Item(45)
$lbl700.Left = EmuToPt 723943
$lbl700.Top = EmuToPt 3908149
$lbl700.Width = EmuToPt 405880
$lbl700.Height = EmuToPt 276999
$lbl700.TextFrame.TextRange.Text = "0%"

# --- "1700" -> "100%" (TextBox 98, id=99) ---
$lbl1700 = $s.Shapes.Item(46)
$lbl1700.Left = EmuToPt 2984528
$lbl1700.Top = EmuToPt 3913802
$lbl1700.Width = EmuToPt 575799
$lbl1700.Height = EmuToPt 276999
$lbl1700.TextFrame.TextRange.Text = "100%"

# --- "1200" -> "75%" (TextBox 102, id=103) ---
$lbl1200 = $s.Shapes.Item(47)
$lbl1200.Left = EmuToPt 2391184
$lbl1200.Top = EmuToPt 3913803
$lbl1200.Width = EmuToPt 490840
$lbl1200.Height = EmuToPt 276999
$lbl1200.TextFrame.TextRange.Text = "75%"

# --- "1000" -> "25%" (TextBox 108, id=109) ---
$lbl1000 = $s.Shapes.Item(48)
$lbl1000.Left = EmuToPt 1348795
$lbl1000.Top = EmuToPt 3901845
$lbl1000.Width = EmuToPt 490840
$lbl1000.Height = EmuToPt 276999
$lbl1000.TextFrame.TextRange.Text = "25%"
